# Daily attendance processing - clear "Recorded By" (column G) free-text
# values for all session rows and shrink the column back down to a normal
# width now that the long recorder-name lists are gone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Recorded By" text for every data row (rows 2-57). Rows 21 and
# 49 were already blank in this column, so clearing them is a no-op.
$ws.Range("G2:G57").ClearContents()

# Column G no longer needs to hold long comma-separated name lists, so
# narrow it from 50 characters down to 13.
$ws.Columns.Item(7).ColumnWidth = 12.17
